$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 2.013740666666667
$ws.Range("H2").Value = 6.041221999999999
$ws.Range("I2").Value = 0.3805515268368102
$ws.Range("J2").Value = 0.3805515268368102
$ws.Range("M2").Value = 4.482602333333333
$ws.Range("N2").Value = 13.447807
$ws.Range("O2").Value = 0.2395001548634358
$ws.Range("P2").Value = 0.2395001548634358
$ws.Range("Q2").Value = 9.026798611128221
$ws.Range("R2").Value = 81.241187500154
$ws.Range("S2").Value = 0.091142149610933
$ws.Range("T2").Value = 0.091142149610933
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 2.013740666666667
$ws.Range("H3").Value = 6.041221999999999
$ws.Range("I3").Value = 0.3805515268368102
$ws.Range("J3").Value = 0.3805515268368102
$ws.Range("O3").Value = 0.1845029314701825
$ws.Range("P3").Value = 0.1845029314701825
$ws.Range("Q3").Value = 6.953944587191556
$ws.Range("R3").Value = 62.585501284724
$ws.Range("S3").Value = 0.07021287227684532
$ws.Range("T3").Value = 0.07021287227684532
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 2.013740666666667
$ws.Range("H4").Value = 6.041221999999999
$ws.Range("I4").Value = 0.3805515268368102
$ws.Range("J4").Value = 0.3805515268368102
$ws.Range("M4").Value = 1.516087333333333
$ws.Range("N4").Value = 4.548262
$ws.Range("O4").Value = 0.08100275779980189
$ws.Range("P4").Value = 0.08100275779980189
$ws.Range("Q4").Value = 3.053006717351555
$ws.Range("R4").Value = 27.477060456164
$ws.Range("S4").Value = 0.03082572315870695
$ws.Range("T4").Value = 0.03082572315870695
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 2.013740666666667
$ws.Range("H5").Value = 6.041221999999999
$ws.Range("I5").Value = 0.3805515268368102
$ws.Range("J5").Value = 0.3805515268368102
$ws.Range("M5").Value = 9.264553333333334
$ws.Range("N5").Value = 27.79366
$ws.Range("O5").Value = 0.4949941558665798
$ws.Range("P5").Value = 0.4949941558665797
$ws.Range("Q5").Value = 18.65640780583556
$ws.Range("R5").Value = 167.90767025252
$ws.Range("S5").Value = 0.188370781790325
$ws.Range("T5").Value = 0.1883707817903249
$ws.Range("G6").Value = 0.8431363333333334
$ws.Range("I6").Value = 0.1593337336295156
$ws.Range("J6").Value = 0.1593337336295156
$ws.Range("M6").Value = 4.482602333333333
$ws.Range("N6").Value = 13.447807
$ws.Range("O6").Value = 0.2395001548634358
$ws.Range("P6").Value = 0.2395001548634358
$ws.Range("Q6").Value = 3.779444895118111
$ws.Range("R6").Value = 34.01500405606301
$ws.Range("S6").Value = 0.03816045387923841
$ws.Range("T6").Value = 0.03816045387923842
$ws.Range("G7").Value = 0.8431363333333334
$ws.Range("I7").Value = 0.1593337336295156
$ws.Range("J7").Value = 0.1593337336295156
$ws.Range("O7").Value = 0.1845029314701825
$ws.Range("P7").Value = 0.1845029314701825
$ws.Range("S7").Value = 0.02939754093673482
$ws.Range("T7").Value = 0.02939754093673483
$ws.Range("G8").Value = 0.8431363333333334
$ws.Range("I8").Value = 0.1593337336295156
$ws.Range("J8").Value = 0.1593337336295156
$ws.Range("M8").Value = 1.516087333333333
$ws.Range("N8").Value = 4.548262
$ws.Range("O8").Value = 0.08100275779980189
$ws.Range("P8").Value = 0.08100275779980189
$ws.Range("Q8").Value = 1.278268315239778
$ws.Range("R8").Value = 11.504414837158
$ws.Range("S8").Value = 0.0129064718345298
$ws.Range("T8").Value = 0.0129064718345298
$ws.Range("G9").Value = 0.8431363333333334
$ws.Range("I9").Value = 0.1593337336295156
$ws.Range("J9").Value = 0.1593337336295156
$ws.Range("M9").Value = 9.264553333333334
$ws.Range("N9").Value = 27.79366
$ws.Range("O9").Value = 0.4949941558665798
$ws.Range("P9").Value = 0.4949941558665797
$ws.Range("Q9").Value = 7.811281527437778
$ws.Range("R9").Value = 70.30153374694001
$ws.Range("S9").Value = 0.07886926697901253
$ws.Range("T9").Value = 0.07886926697901253
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 1.000022
$ws.Range("H10").Value = 3.000066
$ws.Range("I10").Value = 0.1889815830160193
$ws.Range("J10").Value = 0.1889815830160193
$ws.Range("M10").Value = 4.482602333333333
$ws.Range("N10").Value = 13.447807
$ws.Range("O10").Value = 0.2395001548634358
$ws.Range("P10").Value = 0.2395001548634358
$ws.Range("Q10").Value = 4.482700950584667
$ws.Range("R10").Value = 40.34430855526201
$ws.Range("S10").Value = 0.04526111839867387
$ws.Range("T10").Value = 0.04526111839867387
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 1.000022
$ws.Range("H11").Value = 3.000066
$ws.Range("I11").Value = 0.1889815830160193
$ws.Range("J11").Value = 0.1889815830160193
$ws.Range("O11").Value = 0.1845029314701825
$ws.Range("P11").Value = 0.1845029314701825
$ws.Range("Q11").Value = 3.453323304774667
$ws.Range("R11").Value = 31.07990974297201
$ws.Range("S11").Value = 0.03486765606033121
$ws.Range("T11").Value = 0.03486765606033121
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 1.000022
$ws.Range("H12").Value = 3.000066
$ws.Range("I12").Value = 0.1889815830160193
$ws.Range("J12").Value = 0.1889815830160193
$ws.Range("M12").Value = 1.516087333333333
$ws.Range("N12").Value = 4.548262
$ws.Range("O12").Value = 0.08100275779980189
$ws.Range("P12").Value = 0.08100275779980189
$ws.Range("Q12").Value = 1.516120687254667
$ws.Range("R12").Value = 13.645086185292
$ws.Range("S12").Value = 0.01530802939766976
$ws.Range("T12").Value = 0.01530802939766976
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 1.000022
$ws.Range("H13").Value = 3.000066
$ws.Range("I13").Value = 0.1889815830160193
$ws.Range("J13").Value = 0.1889815830160193
$ws.Range("M13").Value = 9.264553333333334
$ws.Range("N13").Value = 27.79366
$ws.Range("O13").Value = 0.4949941558665798
$ws.Range("P13").Value = 0.4949941558665797
$ws.Range("Q13").Value = 9.264757153506668
$ws.Range("R13").Value = 83.38281438156002
$ws.Range("S13").Value = 0.09354477915934443
$ws.Range("T13").Value = 0.09354477915934442
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 1.434738333333333
$ws.Range("H14").Value = 4.304214999999999
$ws.Range("I14").Value = 0.271133156517655
$ws.Range("J14").Value = 0.271133156517655
$ws.Range("M14").Value = 4.482602333333333
$ws.Range("N14").Value = 13.447807
$ws.Range("O14").Value = 0.2395001548634358
$ws.Range("P14").Value = 0.2395001548634358
$ws.Range("Q14").Value = 6.431361400722777
$ws.Range("R14").Value = 57.88225260650499
$ws.Range("S14").Value = 0.06493643297459055
$ws.Range("T14").Value = 0.06493643297459056
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 1.434738333333333
$ws.Range("H15").Value = 4.304214999999999
$ws.Range("I15").Value = 0.271133156517655
$ws.Range("J15").Value = 0.271133156517655
$ws.Range("O15").Value = 0.1845029314701825
$ws.Range("P15").Value = 0.1845029314701825
$ws.Range("Q15").Value = 4.954506323614444
$ws.Range("R15").Value = 44.59055691253
$ws.Range("S15").Value = 0.05002486219627116
$ws.Range("T15").Value = 0.05002486219627117
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 1.434738333333333
$ws.Range("H16").Value = 4.304214999999999
$ws.Range("I16").Value = 0.271133156517655
$ws.Range("J16").Value = 0.271133156517655
$ws.Range("M16").Value = 1.516087333333333
$ws.Range("N16").Value = 4.548262
$ws.Range("O16").Value = 0.08100275779980189
$ws.Range("P16").Value = 0.08100275779980189
$ws.Range("Q16").Value = 2.175188613814444
$ws.Range("R16").Value = 19.57669752433
$ws.Range("S16").Value = 0.02196253340889538
$ws.Range("T16").Value = 0.02196253340889539
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 1.434738333333333
$ws.Range("H17").Value = 4.304214999999999
$ws.Range("I17").Value = 0.271133156517655
$ws.Range("J17").Value = 0.271133156517655
$ws.Range("M17").Value = 9.264553333333334
$ws.Range("N17").Value = 27.79366
$ws.Range("O17").Value = 0.4949941558665798
$ws.Range("P17").Value = 0.4949941558665797
$ws.Range("Q17").Value = 13.29220980854444
$ws.Range("R17").Value = 119.6298882769
$ws.Range("S17").Value = 0.1342093279378979
$ws.Range("T17").Value = 0.1342093279378979
